{"js": "// Change 1: remove the leading \"A use of (quality) software ... [1]. \" sentence\n// from the bullet list item (it stays, but the intro sentence is dropped; the\n// italic sentence that follows it is left untouched).\n{\n  const needle =\n    \"A use of (quality) software helps to produces better scientific discoveries [1]. \";\n  const results = context.document.body.search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].delete();\n    await context.sync();\n  }\n}\n\n// Change 2: reword the \"Role of software in specific domains of research\"\n// Heading 3 into \"Examples of software role in specific domains of research\".\n{\n  const needle = \"Role of software in specific domains of research \";\n  const results = context.document.body.search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Examples of software role in specific domains of research \",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// Change 3: expand the \"Several research breakthroughs...\" sentence into the\n// new combined sentence that also absorbs the \"use of software\" idea that used\n// to live in the bullet list removed above. The following paragraph (\"One of\n// the research breakthroughs is creation of the very first visual\n// representation of a black hole...\") keeps the exact same wording, so it is\n// intentionally left untouched.\n{\n  const needle =\n    \"Several research breakthroughs has been made possible because of the use of software in the research. \";\n  const results = context.document.body.search(needle, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const newText =\n      \"A use of software also allowed to produces better scientific discoveries and several research breakthroughs has been made possible[1]. \";\n    results.items[0].insertText(newText, Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# \"some changes in the role of SW\"\n$d = $word.ActiveDocument\n\n# Change 1: drop the leading \"A use of (quality) software ... [1]. \" sentence\n# from the bullet list item; the italic sentence that follows it\n# (\"A software dictates the quality of a research outcome...\") is left as-is.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"A use of (quality) software helps to produces better scientific discoveries [1]. \"\n$find1.Replacement.Text = \"\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# Change 2: reword the Heading 3 \"Role of software in specific domains of\n# research\" to \"Examples of software role in specific domains of research\".\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"Role of software in specific domains of research \"\n$find2.Replacement.Text = \"Examples of software role in specific domains of research \"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n\n# Change 3: expand \"Several research breakthroughs...\" into the new combined\n# sentence (it now also covers the \"use of software\" idea removed from the\n# bullet list in change 1). The following paragraph, which starts \"One of the\n# research breakthroughs is creation of the very first visual representation\n# of a black hole...\", keeps exactly the same wording, so it is intentionally\n# left untouched.\n$find3 = $d.Content.Find\n$find3.ClearFormatting()\n$find3.Replacement.ClearFormatting()\n$find3.Text = \"Several research breakthroughs has been made possible because of the use of software in the research. \"\n$find3.Replacement.Text = \"A use of software also allowed to produces better scientific discoveries and several research breakthroughs has been made possible[1]. \"\n$find3.Execute($find3.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find3.Replacement.Text, 2) | Out-Null\n"}
